$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("AF", 20, "IDS", 197117),
    @("AN", 22, "IDS", 42150),
    @("CL", 16, "IDS", 117097),
    @("DM", 12, "IDS", 145542),
    @("IA", 23, "IDS", 148972),
    @("JA", 24, "IDS", 47188),
    @("JE", 21, "IDS", 101452),
    @("KA", 14, "IDS", 96879),
    @("KV", 16, "IDS", 106402),
    @("LK", 16, "IDS", 105711),
    @("LM", 12, "IDS", 104149),
    @("LV", 16, "IDS", 148662),
    @("MA", 23, "IDS", 168048),
    @("ME", 12, "IDS", 97909),
    @("MPS", 14, "IDS", 130246),
    @("MS", 23, "IDS", 121924),
    @("MV", 16, "IDS", 194113),
    @("PV", 11, "IDS", 204092),
    @("RV", 11, "IDS", 154644),
    @("SA", 15, "IDS", 97759),
    @("SAF", 19, "IDS", 139226),
    @("SC", 19, "IDS", 171927),
    @("SF", 14, "IDS", 147504),
    @("TL", 11, "IDS", 58618),
    @("TM", 24, "IDS", 101693)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
